$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.781.83"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.06"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.995"
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.92"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  +6.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "25.30"
$ws.Range("E8").Value = "  +8.08%  "
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0590"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0899"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.796.30"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.583.19"
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.784.44"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.95"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.55"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.33"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0689"
$ws.Range("E20").Value = "  +2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.994"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.96"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.12"
$ws.Range("E23").Value = "  +3.68%  "
$ws.Range("E24").Value = "  +4.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.18"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("E26").Value = "  +4.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.88"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.30"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0462"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.412.89"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.02"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.05"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.75"
$ws.Range("E37").Value = "  +6.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.29"
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0163"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.522"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.95"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.996"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.779"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.95"
$ws.Range("E45").Value = "  +3.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.27"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.710.04"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.98"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.828"
$ws.Range("E49").Value = "  -8.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.50"
$ws.Range("E50").Value = "  +2.97%  "
$ws.Range("E51").Value = "  +0.31%  "
